$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.925.27'
$ws.Range('E2').Value = '  -3.76%  '
$ws.Range('D3').Value = '1.636.68'
$ws.Range('E3').Value = '  -6.04%  '
$ws.Range('D4').Value = '''0.9976'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''235.86'
$ws.Range('E5').Value = '  -4.58%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '''0.4691'
$ws.Range('E7').Value = '  -7.17%  '
$ws.Range('D8').Value = '''0.2557'
$ws.Range('E8').Value = '  -5.80%  '
$ws.Range('D9').Value = '''0.06007'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('D10').Value = '''0.07124'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').Value = '1.636.74'
$ws.Range('E11').Value = '  -6.02%  '
$ws.Range('D12').Value = '''14.80'
$ws.Range('E12').Value = '  -1.88%  '
$ws.Range('D13').Value = '''0.6144'
$ws.Range('E13').Value = '  -5.14%  '
$ws.Range('D14').Value = '''4.396'
$ws.Range('E14').Value = '  -4.97%  '
$ws.Range('D15').Value = '''72.60'
$ws.Range('E15').Value = '  -6.31%  '
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '''0.9978'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '24.924.34'
$ws.Range('E18').Value = '  -3.86%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '''0.000006575'
$ws.Range('E19').Value = '  -3.32%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '''11.23'
$ws.Range('E20').Value = '  -4.92%  '
$ws.Range('D21').Value = '''4.397'
$ws.Range('E21').Value = '  +2.68%  '
$ws.Range('D22').Value = '1.842.93'
$ws.Range('E22').Value = '  -6.37%  '
$ws.Range('D23').Value = '''8.559'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('D24').Value = '''5.248'
$ws.Range('E24').Value = '  -2.38%  '
$ws.Range('D25').Value = '''132.26'
$ws.Range('E25').Value = '  -2.89%  '
$ws.Range('D27').Value = '''1.368'
$ws.Range('E27').Value = '  -8.89%  '
$ws.Range('D28').Value = '''102.45'
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('D29').Value = '''1.653'
$ws.Range('E29').Value = '  -6.34%  '
$ws.Range('D30').Value = '''3.723'
$ws.Range('E30').Value = '  -4.77%  '
$ws.Range('D31').Value = '''0.07730'
$ws.Range('E31').Value = '  -6.03%  '
$ws.Range('D32').Value = '''3.538'
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('D33').Value = '''0.04363'
$ws.Range('E33').Value = '  -6.48%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').Value = '''2.594'
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('D36').Value = '''0.9187'
$ws.Range('E36').Value = '  -7.45%  '
$ws.Range('D37').Value = '''0.5799'
$ws.Range('E37').Value = '  -6.36%  '
$ws.Range('D38').Value = '''2.544'
$ws.Range('E38').Value = '  -6.42%  '
$ws.Range('D39').Value = '''0.01552'
$ws.Range('E39').Value = '  -3.05%  '
$ws.Range('D40').Value = '''0.9977'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '''0.8121'
$ws.Range('E41').Value = '  +7.28%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''1.795'
$ws.Range('E42').Value = '  -6.02%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''97.38'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.3707'
$ws.Range('E44').Value = '  -3.78%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''4.729'
$ws.Range('E45').Value = '  -5.05%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '''0.1121'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.05224'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''6.081'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''29.51'
$ws.Range('E49').Value = '  -3.55%  '
$ws.Range('B50').Value = 'TrueUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D50').Value = '''0.9995'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '''0.9993'
$ws.Range('E51').Value = '  -0.49%  '
